# Quarterly database update (rolling one quarter forward) + EPS restatement
# for the "1400/09" quarter (read_price algorithm change).
#
# Every quarter column D..M is shifted one step to the left (the oldest
# quarter in D drops off) and the newly-reported quarter's figures are
# appended in column M. Mirrors the author's "update database and change
# read_price algorithm" commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Row 8: quarter-name header -------------------------------------------------
$ws.Range("D8").Value = "فصل سوم منتهی به 1399/09"
$ws.Range("E8").Value = "فصل چهارم منتهی به 1399/12"
$ws.Range("F8").Value = "فصل اول منتهی به 1400/03"
$ws.Range("G8").Value = "فصل دوم منتهی به 1400/06"
$ws.Range("H8").Value = "فصل سوم منتهی به 1400/09"
$ws.Range("I8").Value = "فصل چهارم منتهی به 1400/12"
$ws.Range("J8").Value = "فصل اول منتهی به 1401/03"
$ws.Range("K8").Value = "فصل دوم منتهی به 1401/06"
$ws.Range("L8").Value = "فصل سوم منتهی به 1401/09"
$ws.Range("M8").Value = "فصل چهارم منتهی به 1401/12"

# ---- Row 9: publish-date header ------------------------------------------------
$ws.Range("D9").Value = "1400-10-29 (2)"
$ws.Range("E9").Value = "1401-04-21 (10)"
$ws.Range("F9").Value = "1401-04-30 (2)"
$ws.Range("G9").Value = "1401-09-23 (6)"
$ws.Range("H9").Value = "1401-10-28 (2)"
$ws.Range("I9").Value = "1402-02-29 (8)"
$ws.Range("J9").Value = "1401-04-30"
$ws.Range("K9").Value = "1401-09-23 (3)"
$ws.Range("L9").Value = "1401-10-28"
$ws.Range("M9").Value = "1402-02-29"

# ---- Data rows: roll each quarter left by one column, append new quarter ------
# row -> ordered values for columns D,E,F,G,H,I,J,K,L,M
$rows = @{
    11 = @(1288603, 1729870, 1781465, 1668372, 2011673, 2566637, 3201905, 2690382, 2183035, 1899853)
    12 = @(-1105859, -1477904, -1639142, -1568793, -1829511, -2265612, -2597072, -2436766, -1985710, -1636715)
    13 = @(182744, 251966, 142323, 99580, 182163, 301026, 604833, 253615, 197324, 263138)
    14 = @(-8455, -10206, -13025, -13967, -12932, -8178, -27388, -8820, -17604, -20395)
    16 = @(535, 19050, "-", "-", 43148, -1090, -208, "-", "-", -353)
    17 = @(174824, 260810, 129297, 85613, 212378, 291758, 577237, 244795, 179721, 242389)
    18 = @("-", -363, -380, -7219, -8136, -8190, -12945, -13815, -11053, -9863)
    19 = @(20279, -10218, 6404, 12148, 16011, -4220, -6204, 36637, 1598, 34735)
    20 = @(195103, 250228, 135321, 90542, 220254, 279348, 558088, 267618, 170265, 267262)
    21 = @(-38164, 5120, -25514, -15589, -41734, -26762, -111206, -9877, -24506, -3330)
    22 = @(156939, 255348, 109807, 74953, 178520, 252585, 446882, 257742, 145759, 263931)
    24 = @(156939, 255348, 109807, 74953, 178520, 252585, 446882, 257742, 145759, 263931)
    26 = @(279659, 567669, 597835, 727938, 667334, 689812, 646594, 612590, 546874, 499416)
}

$cols = @("D", "E", "F", "G", "H", "I", "J", "K", "L", "M")

foreach ($r in $rows.Keys) {
    $vals = $rows[$r]
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $ws.Range($cols[$i] + $r).Value = $vals[$i]
    }
}
